# Overwrite old files with RMI version
# - "About" sheet: drop the old citation block (date stamp, author, title,
#   hyperlink, page reference) and replace the "Source:" value with "None";
#   shift the "Notes" block up and append a new note explaining that the US
#   version zeroes out this variable.
# - "EoCSoEVMS" sheet: zero out the market-share-change coefficient so
#   additional chargers no longer induce extra EV adoption in the US model.

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$data  = $wb.Worksheets.Item("EoCSoEVMS")

# --- Remove the old hyperlink on the About sheet before we clear its cell ---
foreach ($hl in $about.Hyperlinks) {
    $hl.Delete()
}

# --- Wipe the old "About" content (values + formatting) clean so the old
#     citation rows don't leave stray formatted-but-empty cells behind ---
$about.UsedRange.Clear()

# --- Rebuild the About sheet content per the new layout ---
$about.Range("A1").Value = "EoCSoEVMS Effect of Charging Stations on EV Market Share"
$about.Range("A1").Font.Bold = $true

$about.Range("A3").Value = "Source:"
$about.Range("A3").Font.Bold = $true

$about.Range("A5").Value = "Notes"
$about.Range("A5").Font.Bold = $true

$about.Range("A6").Value = "This variable must be expressed as the percent increase in EV market share"
$about.Range("A7").Value = "for every 1 unit increase in ""EV chargers per 100,000 people""."

$about.Range("A9").Value = "In the US, we set this to 0 so that increasing EV chargers does not induce additional deployment."

$about.Range("B3").Value = "None"

# --- The "Hyperlink" cell style is no longer used anywhere in the workbook
#     now that the citation hyperlink is gone; drop it so styles.xml doesn't
#     keep carrying a dead named style around. ---
foreach ($st in $wb.Styles) {
    if ($st.Name() -eq "Hyperlink") {
        $st.Delete()
    }
}

# --- Update the EoCSoEVMS data sheet: zero out the coefficient ---
$data.Range("B2").Value = 0
